{"js": "// The document contains three <id>...</id> tags rendered across separate\n// runs, e.g. \"<id>\" + \"p092v_a1\" + \"</id>\". This edit collapses each of\n// them into a single run of text and renumbers the inner id from\n// \"p092v_aN\" to \"p092v_N\" (dropping the \"a\"):\n//   <id>p092v_a1</id> -> <id>p092v_1</id>\n//   <id>p092v_a2</id> -> <id>p092v_2</id>\n//   <id>p092v_a3</id> -> <id>p092v_3</id>\n\nconst renumber = [\n  [\"p092v_a1\", \"p092v_1\"],\n  [\"p092v_a2\", \"p092v_2\"],\n  [\"p092v_a3\", \"p092v_3\"],\n];\n\nfor (const [oldId, newId] of renumber) {\n  const oldText = `<id>${oldId}</id>`;\n  const newText = `<id>${newId}</id>`;\n\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    // Replacing the whole matched range (which spans the original\n    // \"<id>\", \"p092v_aN\" and \"</id>\" runs) collapses it into a single\n    // run that keeps the formatting of the range's first run.\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains three <id>...</id> tags rendered across separate\n# runs, e.g. \"<id>\" + \"p092v_a1\" + \"</id>\". This edit collapses each of\n# them into a single run of text and renumbers the inner id from\n# \"p092v_aN\" to \"p092v_N\" (dropping the \"a\"):\n#   <id>p092v_a1</id> -> <id>p092v_1</id>\n#   <id>p092v_a2</id> -> <id>p092v_2</id>\n#   <id>p092v_a3</id> -> <id>p092v_3</id>\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"p092v_a1\", \"p092v_1\"),\n  @(\"p092v_a2\", \"p092v_2\"),\n  @(\"p092v_a3\", \"p092v_3\")\n)\n\nforeach ($pair in $pairs) {\n  $oldId = $pair[0]\n  $newId = $pair[1]\n  $oldText = \"<id>$oldId</id>\"\n  $newText = \"<id>$newId</id>\"\n\n  # Search the whole tag (opening + inner id + closing) as one string so\n  # the match spans all three original runs; replacing it collapses the\n  # match into a single run that keeps the formatting of the first run\n  # in the match (the Courier New \"<id>\" run).\n  $range = $d.Content\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
